# Add a primary (default) header to the document containing the
# questionnaire number, so printed copies can be tracked.
$d = $word.ActiveDocument

# wdHeaderFooterPrimary = 1
$section = $d.Sections.First
$header = $section.Headers.Item(1)

# Insert the header text. Using InsertAfter (rather than assigning to
# .Range.Text) keeps this a simple single "default" header instead of
# also materializing even/first header & footer variants.
$header.Range.InsertAfter("Questionnaire 3")

# Apply the "Header" paragraph style and center it.
$header.Range.Style = "Header"
$header.Range.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter

# Format just the inserted run (exclude the trailing paragraph mark) as
# centered 12pt Arial.
$textRange = $header.Range.Duplicate
$textRange.End = $textRange.End - 1
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12

Write-Output "Header added."
